# Reorder worksheets: "Статистика" moves to the front, "Ошибки заполнения"
# stays in the middle, "Модель" moves to the end. Also update the active/
# selected cells on a couple of sheets and nudge the chart picture on the
# "Статистика" sheet down a bit.

$wb = $excel.ActiveWorkbook

# --- Reorder the sheet tabs -------------------------------------------------
$stat = $wb.Worksheets.Item("Статистика")
$stat.Move($wb.Worksheets.Item(1))

$model = $wb.Worksheets.Item("Модель")
$model.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# --- Make "Статистика" the active tab and update its selection -------------
$statSheet = $wb.Worksheets.Item("Статистика")
$statSheet.Activate()
$statSheet.Range("I5").Select()

# --- Move the chart picture on the "Статистика" sheet down -----------------
$plot = $statSheet.Shapes.Item("MyPlot")
$plot.Top = 86.4

# --- Reset the selection on "Модель" and scroll it over to column F --------
$modelSheet = $wb.Worksheets.Item("Модель")
$modelSheet.Activate()
$modelSheet.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1

# --- Re-activate "Статистика" so it is the tab shown when the file opens ---
$statSheet.Activate()
